$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "94.930.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.485.48"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.49%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "645.66"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.47"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +4.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.407"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.55%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.483.74"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.63"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.63%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.140.16"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.75%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "94.733.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000257"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.58%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.471.37"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.02"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.02%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.53"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +9.43%  "

$ws.Range("B22").Value = "Stellar"
$ws.Range("C22").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.518"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -8.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "505.76"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000194"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.68"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.38"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.18"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.672.67"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +4.77%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.93"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +7.52%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +11.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.139"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.50%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.07"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +10.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.574"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.82"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "544.36"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.945"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +13.75%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.72"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.71"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0417"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.55"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.17"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +9.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.23"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.37"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.15%  "
